# Applies the two classes of changes from the commit:
# 1. Column C ("Förändrad") date value changes from 45184 to 45186 for every
#    data row (rows 2-271).
# 2. Every HYPERLINK(...) formula in columns S-Y gets a second argument added:
#    the "friendly name" to display, which equals the row's "Beteckning"
#    (column A) value - the same text embedded in the URL (filename, no ext).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 271

# --- 1. Update column C date values for all data rows ---
$dateRangeAddr = "C" + $firstRow + ":C" + $lastRow
$ws.Range($dateRangeAddr).Value = 45186

# --- 2. Update HYPERLINK formulas to include the friendly display name ---
$hyperlinkCols = @("S","T","U","V","W","X","Y")

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $name = $ws.Range("A$r").Value2
    if ([string]::IsNullOrEmpty($name)) {
        continue
    }

    foreach ($col in $hyperlinkCols) {
        $addr = "$col$r"
        $cell = $ws.Range($addr)
        $f = $cell.Formula

        if ($f -ne "" -and $f -like "*HYPERLINK(*" -and $f -notlike "*,*") {
            # Insert the friendly name as a second argument right before the
            # formula's final closing parenthesis.
            $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $name + '")'
            $cell.Formula = $newFormula
        }
    }
}
